$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 16 (total_venda for day 21/07 changed).
$ws.Cells.Item(16, 2).Value2 = 12417.3

# A new daily entry (day 22, July/2025) was added, inserted as a new row 17,
# pushing every row below it down by one.
$ws.Rows.Item(17).Insert()

$ws.Cells.Item(17, 1).Value2 = 22
$ws.Cells.Item(17, 2).Value2 = 7698.86
$ws.Cells.Item(17, 3).Value2 = 7
$ws.Cells.Item(17, 4).Value2 = 2025
$ws.Cells.Item(17, 5).Value2 = "07/2025"
